$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update individual G-column (Waargenomen) corrections for existing weeks ---
$ws.Range("G10").Value = 3379
$ws.Range("G14").Value = 2727
$ws.Range("G19").Value = 2638
$ws.Range("G22").Value = 2671
$ws.Range("G23").Value = 2665
$ws.Range("G24").Value = 2639
$ws.Range("G26").Value = 2850
$ws.Range("G28").Value = 2688
$ws.Range("G30").Value = 2717
$ws.Range("G31").Value = 2889
$ws.Range("G32").Value = 2996
$ws.Range("G35").Value = 3441
$ws.Range("G36").Value = 3670
$ws.Range("G37").Value = 3574
$ws.Range("G38").Value = 3530

# --- Add new week 47 row (row 39) ---
$ws.Range("F39").Value = 47
$ws.Range("G39").Value = 3336
$ws.Range("H39").Value = 2972
$ws.Range("I39").Formula = "=G39-H39"

# --- Move the "Som week 11 tot en met 19" total row from row 40 to row 42 ---
$ws.Range("F42").Value = "Som week 11 tot en met 19"
$ws.Range("G42").Formula = "=SUM(G3:G28)"
$ws.Range("H42").Formula = "=SUM(H3:H28)"
$ws.Range("I42").Formula = "=SUM(I3:I34)"
$ws.Range("G42:I42").NumberFormat = "0"

$ws.Range("F40:I40").Clear()

# --- Restore the view state recorded in the saved workbook ---
$ws.Range("L37").Select()

$wb.Save()
